$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 21:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1633639
$ws.Range("C4").Value = 12737
$ws.Range("D4").Value = 385286
$ws.Range("E4").Value = 1151184
$ws.Range("G4").Value = 815
$ws.Range("H4").Value = 97169

# Row 10 - Francia
$ws.Range("D10").Value = 64209
$ws.Range("E10").Value = 89328
$ws.Range("G10").Value = 74
$ws.Range("H10").Value = 28289

# Row 11 - Alemania
$ws.Range("B11").Value = 179626
$ws.Range("C11").Value = 605
$ws.Range("E11").Value = 12282
$ws.Range("G11").Value = 35
$ws.Range("H11").Value = 8344

# Row 14 - India
$ws.Range("B14").Value = 124747
$ws.Range("C14").Value = 6521
$ws.Range("D14").Value = 51807
$ws.Range("E14").Value = 69214
$ws.Range("G14").Value = 142
$ws.Range("H14").Value = 3726

# Row 101 - Maldivas
$ws.Range("B101").Value = 1274
$ws.Range("C101").Value = 58
$ws.Range("E101").Value = 1172

# Row 115 - Costa Rica
$ws.Range("B115").Value = 911
$ws.Range("C115").Value = 8
$ws.Range("D115").Value = 600

$wb.Save()
